$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("7:7").Insert()
Write-Host "insert ok"
